$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Case_3_17 (380 kV) loading_percent.xlsx results: update computed loading
# percentages for rows 2-25 (time steps 0-23) across line columns
# B,C,D,F,G,J,K,M,O (the zero columns E,H,I,L,N are unchanged).

$ws.Range("B2").Value = 9.873373792187643
$ws.Range("C2").Value = 4.535836034042181
$ws.Range("D2").Value = 9.872111394706771
$ws.Range("F2").Value = 36.12361155250988
$ws.Range("G2").Value = 3.67433603405399
$ws.Range("J2").Value = 11.46120904569205
$ws.Range("K2").Value = 9.236550822851603
$ws.Range("M2").Value = 15.97307109453593
$ws.Range("O2").Value = 26.9955597954439
$ws.Range("B3").Value = 9.630720560333298
$ws.Range("C3").Value = 4.402658554763506
$ws.Range("D3").Value = 9.82910144766784
$ws.Range("F3").Value = 36.18940175905835
$ws.Range("G3").Value = 3.6762383267413
$ws.Range("J3").Value = 11.47997465097203
$ws.Range("K3").Value = 9.076774920515803
$ws.Range("M3").Value = 15.90293171799192
$ws.Range("O3").Value = 27.07714197310857
$ws.Range("B4").Value = 9.480387537821546
$ws.Range("C4").Value = 4.318145172239533
$ws.Range("D4").Value = 9.804321696640667
$ws.Range("F4").Value = 36.237692912373
$ws.Range("G4").Value = 3.677468406910593
$ws.Range("J4").Value = 11.49313265410786
$ws.Range("K4").Value = 8.978686574455431
$ws.Range("M4").Value = 15.86237874199631
$ws.Range("O4").Value = 27.13245152654292
$ws.Range("B5").Value = 9.41888126546386
$ws.Range("C5").Value = 4.283049388703776
$ws.Range("D5").Value = 9.794640898560795
$ws.Range("F5").Value = 36.25935398026564
$ws.Range("G5").Value = 3.6779853302371
$ws.Range("J5").Value = 11.49890605415958
$ws.Range("K5").Value = 8.938769453142564
$ws.Range("M5").Value = 15.84649736969222
$ws.Range("O5").Value = 27.15630056274017
$ws.Range("B6").Value = 9.408656244067942
$ws.Range("C6").Value = 4.277183184583472
$ws.Range("D6").Value = 9.793058817946671
$ws.Range("F6").Value = 36.26307038439414
$ws.Range("G6").Value = 3.678072111967171
$ws.Range("J6").Value = 11.49988957173186
$ws.Range("K6").Value = 8.932146050503373
$ws.Range("M6").Value = 15.84389954146002
$ws.Range("O6").Value = 27.16033973128692
$ws.Range("B7").Value = 9.479558905774862
$ws.Range("C7").Value = 4.317674467416856
$ws.Range("D7").Value = 9.804189439247409
$ws.Range("F7").Value = 36.2379770208463
$ws.Range("G7").Value = 3.677475314865549
$ws.Range("J7").Value = 11.49320885038135
$ws.Range("K7").Value = 8.978147948974341
$ws.Range("M7").Value = 15.86216193570179
$ws.Range("O7").Value = 27.13276786193737
$ws.Range("B8").Value = 9.79004041192383
$ws.Range("C8").Value = 4.490504753953067
$ws.Range("D8").Value = 9.856948296620596
$ws.Range("F8").Value = 36.14465535948047
$ws.Range("G8").Value = 3.674979092563009
$ws.Range("J8").Value = 11.46733995433732
$ws.Range("K8").Value = 9.181487391112324
$ws.Range("M8").Value = 15.94837315603722
$ws.Range("O8").Value = 27.02260514175205
$ws.Range("B9").Value = 10.38435737255254
$ws.Range("C9").Value = 4.806341561348219
$ws.Range("D9").Value = 9.972972547315093
$ws.Range("F9").Value = 36.02443426906659
$ws.Range("G9").Value = 3.670574227980445
$ws.Range("J9").Value = 11.42958731620401
$ws.Range("K9").Value = 9.578157262088444
$ws.Range("M9").Value = 16.1367940782654
$ws.Range("O9").Value = 26.84806853190904
$ws.Range("B10").Value = 10.80728998409366
$ws.Range("C10").Value = 5.022736700796913
$ws.Range("D10").Value = 10.06536705441441
$ws.Range("F10").Value = 35.97453603175511
$ws.Range("G10").Value = 3.667633692846355
$ws.Range("J10").Value = 11.40975499924327
$ws.Range("K10").Value = 9.865472897447068
$ws.Range("M10").Value = 16.28621200761846
$ws.Range("O10").Value = 26.74525233709999
$ws.Range("B11").Value = 10.99579629627072
$ws.Range("C11").Value = 5.117501878347336
$ws.Range("D11").Value = 10.10883580626156
$ws.Range("F11").Value = 35.96020183889172
$ws.Range("G11").Value = 3.666359513610046
$ws.Range("J11").Value = 11.40244695962977
$ws.Range("K11").Value = 9.994714882707525
$ws.Range("M11").Value = 16.35638428073993
$ws.Range("O11").Value = 26.70402151779414
$ws.Range("B12").Value = 11.06655064475224
$ws.Range("C12").Value = 5.152838513959709
$ws.Range("D12").Value = 10.12549265919067
$ws.Range("F12").Value = 35.95597735932286
$ws.Range("G12").Value = 3.665886093716519
$ws.Range("J12").Value = 11.39992574907076
$ws.Range("K12").Value = 10.04340065069156
$ws.Range("M12").Value = 16.38325659774022
$ws.Range("O12").Value = 26.6892070703177
$ws.Range("B13").Value = 11.05134145217559
$ws.Range("C13").Value = 5.145252851960167
$ws.Range("D13").Value = 10.12189675485669
$ws.Range("F13").Value = 35.95683364286668
$ws.Range("G13").Value = 3.66598764985401
$ws.Range("J13").Value = 11.40045779235468
$ws.Range("K13").Value = 10.03292735757153
$ws.Range("M13").Value = 16.37745611331802
$ws.Range("O13").Value = 26.69236207044458
$ws.Range("B14").Value = 11.00163028309934
$ws.Range("C14").Value = 5.120420154649394
$ws.Range("D14").Value = 10.11020230567946
$ws.Range("F14").Value = 35.95983016525068
$ws.Range("G14").Value = 3.666320383266941
$ws.Range("J14").Value = 11.40223460585149
$ws.Range("K14").Value = 9.998725660735694
$ws.Range("M14").Value = 16.35858917668302
$ws.Range("O14").Value = 26.70278670467616
$ws.Range("B15").Value = 10.9710968792676
$ws.Range("C15").Value = 5.105137359016056
$ws.Range("D15").Value = 10.10306434626258
$ws.Range("F15").Value = 35.96182237238352
$ws.Range("G15").Value = 3.666525373784432
$ws.Range("J15").Value = 11.40335500753732
$ws.Range("K15").Value = 9.97774156526102
$ws.Range("M15").Value = 16.34707113924437
$ws.Range("O15").Value = 26.70927617506976
$ws.Range("B16").Value = 10.79488632092298
$ws.Range("C16").Value = 5.016467832394519
$ws.Range("D16").Value = 10.06255430654955
$ws.Range("F16").Value = 35.97564120205563
$ws.Range("G16").Value = 3.667718237086341
$ws.Range("J16").Value = 11.41026706201674
$ws.Range("K16").Value = 9.856993322193919
$ws.Range("M16").Value = 16.28166894267414
$ws.Range("O16").Value = 26.74805853748086
$ws.Range("B17").Value = 10.68573959285112
$ws.Range("C17").Value = 4.961115573426945
$ws.Range("D17").Value = 10.03806357760893
$ws.Range("F17").Value = 35.9862617425352
$ws.Range("G17").Value = 3.668466248020087
$ws.Range("J17").Value = 11.4149461607511
$ws.Range("K17").Value = 9.782511324372624
$ws.Range("M17").Value = 16.24209871623604
$ws.Range("O17").Value = 26.77327085693717
$ws.Range("B18").Value = 10.62259915271745
$ws.Range("C18").Value = 4.928933675116862
$ws.Range("D18").Value = 10.02411316811194
$ws.Range("F18").Value = 35.99315772463944
$ws.Range("G18").Value = 3.66890246229158
$ws.Range("J18").Value = 11.41779878244538
$ws.Range("K18").Value = 9.739536474613301
$ws.Range("M18").Value = 16.21954721205775
$ws.Range("O18").Value = 26.78829372805895
$ws.Range("B19").Value = 10.60116090923932
$ws.Range("C19").Value = 4.917978891594967
$ws.Range("D19").Value = 10.01941347698017
$ws.Range("F19").Value = 35.99562777751559
$ws.Range("G19").Value = 3.669051185079473
$ws.Range("J19").Value = 11.41879234618727
$ws.Range("K19").Value = 9.724964204649224
$ws.Range("M19").Value = 16.21194793480189
$ws.Range("O19").Value = 26.79346970387454
$ws.Range("B20").Value = 10.69739645118614
$ws.Range("C20").Value = 4.967043751674828
$ws.Range("D20").Value = 10.04065665682692
$ws.Range("F20").Value = 35.98504967505006
$ws.Range("G20").Value = 3.668386002588607
$ws.Range("J20").Value = 11.41443136753878
$ws.Range("K20").Value = 9.790454345684051
$ws.Range("M20").Value = 16.24628960323802
$ws.Range("O20").Value = 26.77053298335253
$ws.Range("B21").Value = 11.01624925514768
$ws.Range("C21").Value = 5.127729168781557
$ws.Range("D21").Value = 10.11363201066055
$ws.Range("F21").Value = 35.95891734690488
$ws.Range("G21").Value = 3.666222405213186
$ws.Range("J21").Value = 11.40170603355769
$ws.Range("K21").Value = 10.0087788102554
$ws.Range("M21").Value = 16.36412285956177
$ws.Range("O21").Value = 26.69970304311629
$ws.Range("B22").Value = 11.22094385870175
$ws.Range("C22").Value = 5.229539121740863
$ws.Range("D22").Value = 10.16246398292122
$ws.Range("F22").Value = 35.94885355578333
$ws.Range("G22").Value = 3.664861297927012
$ws.Range("J22").Value = 11.39482409168206
$ws.Range("K22").Value = 10.14995850422524
$ws.Range("M22").Value = 16.44287183466853
$ws.Range("O22").Value = 26.65806826591177
$ws.Range("B23").Value = 11.11205412344884
$ws.Range("C23").Value = 5.175500907383634
$ws.Range("D23").Value = 10.13630086147842
$ws.Range("F23").Value = 35.95358282691669
$ws.Range("G23").Value = 3.665582918275623
$ws.Range("J23").Value = 11.39836592779653
$ws.Range("K23").Value = 10.0747606699781
$ws.Range("M23").Value = 16.40068869844967
$ws.Range("O23").Value = 26.67986281037686
$ws.Range("B24").Value = 10.69212759658564
$ws.Range("C24").Value = 4.964364739391631
$ws.Range("D24").Value = 10.03948392111506
$ws.Range("F24").Value = 35.9855951897106
$ws.Range("G24").Value = 3.668422262317058
$ws.Range("J24").Value = 11.41466359921639
$ws.Range("K24").Value = 9.786863782424865
$ws.Range("M24").Value = 16.24439428536192
$ws.Range("O24").Value = 26.77176913152005
$ws.Range("B25").Value = 10.22565633456206
$ws.Range("C25").Value = 4.723544023047913
$ws.Range("D25").Value = 9.940290548344736
$ws.Range("F25").Value = 36.0502177608055
$ws.Range("G25").Value = 3.671713704284216
$ws.Range("J25").Value = 11.43841164621968
$ws.Range("K25").Value = 9.471370387055767
$ws.Range("M25").Value = 16.08383033679718
$ws.Range("O25").Value = 26.89083039504157

Write-Output "Updated $(${ws}.UsedRange.Cells.Count) cells in loading_percent sheet"
